$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New steady-state voltage results (p.u.) for the 380 kV case (Case_0_243)
# Maps row number -> column letter -> updated value, matching res_bus/vm_pu.xlsx
$updates = @{
    2 = @{ "B"=1.02; "C"=1.023668524355633; "D"=1.033298432064749; "E"=1.024233148787848; "F"=1.041205138650722; "I"=1.026690890226896; "J"=1.028847451183388; "K"=1.03610118353107; "L"=1.027062251678982; "M"=1.043985313303571; "N"=1.013584193440138 }
    3 = @{ "B"=1.02; "C"=1.024805691979436; "D"=1.034374192760211; "E"=1.025204607911374; "F"=1.042348178194984; "I"=1.026662920049573; "J"=1.029622223995774; "K"=1.036985183061081; "L"=1.027840311031423; "M"=1.044938059856766; "N"=1.013848126944308 }
    4 = @{ "B"=1.02; "C"=1.025541807101835; "D"=1.035070700469297; "E"=1.025833788668838; "F"=1.043087701126899; "I"=1.026642424156548; "J"=1.030123340533146; "K"=1.037557017544349; "L"=1.028343737488784; "M"=1.045553872913479; "N"=1.014018644641665 }
    5 = @{ "B"=1.02; "C"=1.025851340542747; "D"=1.035363612640163; "E"=1.026098435507715; "F"=1.043398572538337; "I"=1.026633232787293; "J"=1.030333959015028; "K"=1.037797375573707; "L"=1.028555370856397; "M"=1.045812598255217; "N"=1.014090266791144 }
    6 = @{ "B"=1.02; "C"=1.025903316736159; "D"=1.035412799769255; "E"=1.026142879041402; "F"=1.04345076780507; "I"=1.02663165577525; "J"=1.03036931980382; "K"=1.037837730331142; "L"=1.028590904614558; "M"=1.045856029844551; "N"=1.014102288740146 }
    7 = @{ "B"=1.02; "C"=1.02554594282421; "D"=1.035074613981707; "E"=1.025837324345739; "F"=1.043091855100463; "I"=1.026642303601218; "J"=1.03012615502748; "K"=1.037560229381878; "L"=1.028346565372095; "M"=1.045557330652883; "N"=1.014019601909893 }
    8 = @{ "B"=1.02; "C"=1.024052775591263; "D"=1.033661903823305; "E"=1.024561337479824; "F"=1.041591454900157; "I"=1.026681933373675; "J"=1.029109333533941; "K"=1.036399971033078; "L"=1.027325206722249; "M"=1.044307438640736; "N"=1.013673445819455 }
    9 = @{ "B"=1.02; "C"=1.021423820390545; "D"=1.03117572823494; "E"=1.022317341476179; "F"=1.038946793612392; "I"=1.026733447910888; "J"=1.027315923508231; "K"=1.034354124138841; "L"=1.025525209549109; "M"=1.04209978844082; "N"=1.013061448223862 }
    10 = @{ "B"=1.02; "C"=1.019672605538829; "D"=1.029520418090541; "E"=1.020824337628465; "F"=1.037183170966654; "I"=1.026755531297127; "J"=1.026119201232412; "K"=1.032989328321303; "L"=1.024325047967675; "M"=1.040624547530097; "N"=1.012652090106604 }
    11 = @{ "B"=1.02; "C"=1.018914634472028; "D"=1.028804154604117; "E"=1.020178559035151; "F"=1.036419379170392; "I"=1.026762197512319; "J"=1.025600739100134; "K"=1.032398140962829; "L"=1.023805323634868; "M"=1.039984928713302; "N"=1.012474511263634 }
    12 = @{ "B"=1.02; "C"=1.01863313660767; "D"=1.02853817666281; "E"=1.019938793343725; "F"=1.036135652849299; "I"=1.026764239446237; "J"=1.025408117690654; "K"=1.032178513994378; "L"=1.02361226760583; "M"=1.039747221068975; "N"=1.012408501866185 }
    13 = @{ "B"=1.02; "C"=1.018693516777726; "D"=1.028595226482732; "E"=1.019990219140128; "F"=1.036196514019644; "I"=1.026763821078379; "J"=1.025449437520764; "K"=1.03222562625342; "L"=1.023653679113154; "M"=1.039798215811438; "N"=1.0124226633105 }
    14 = @{ "B"=1.02; "C"=1.018891364832479; "D"=1.028782167273104; "E"=1.020158737782921; "F"=1.036395926674612; "I"=1.026762375151495; "J"=1.025584817803023; "K"=1.032379987217774; "L"=1.023789365713313; "M"=1.03996528227031; "N"=1.012469055900005 }
    15 = @{ "B"=1.02; "C"=1.0190132716528; "D"=1.028897357484743; "E"=1.02026258161799; "F"=1.036518788728925; "I"=1.026761426759896; "J"=1.025668224553338; "K"=1.032475089631751; "L"=1.023872965733616; "M"=1.040068200898373; "N"=1.012497633450596 }
    16 = @{ "B"=1.02; "C"=1.019722916147193; "D"=1.029567964569007; "E"=1.020867210596486; "F"=1.037233858517602; "I"=1.026755027943447; "J"=1.02615360406037; "K"=1.033028558777702; "L"=1.024359539356406; "M"=1.040666979436045; "N"=1.012663868599999 }
    17 = @{ "B"=1.02; "C"=1.020168141273774; "D"=1.029988751369867; "E"=1.021246666148032; "F"=1.037682367491234; "I"=1.026750239424003; "J"=1.026457996298825; "K"=1.033375675736895; "L"=1.024664741493466; "M"=1.041042355191258; "N"=1.012768056739419 }
    18 = @{ "B"=1.02; "C"=1.020427864178542; "D"=1.030234237147057; "E"=1.021468063892147; "F"=1.037943962513196; "I"=1.026747166876285; "J"=1.026635516780729; "K"=1.033578121812653; "L"=1.024842756274537; "M"=1.041261225463077; "N"=1.012828796651083 }
    19 = @{ "B"=1.02; "C"=1.020516428208068; "D"=1.030317949599936; "E"=1.021543566298284; "F"=1.038033157514037; "I"=1.026746071783743; "J"=1.026696042218418; "K"=1.033647147077869; "L"=1.024903453966647; "M"=1.041335840993055; "N"=1.012849502088397 }
    20 = @{ "B"=1.02; "C"=1.020120369672549; "D"=1.029943599971018; "E"=1.02120594715322; "F"=1.037634248035736; "I"=1.026750782088746; "J"=1.026425340621756; "K"=1.033338435585365; "L"=1.024631996656058; "M"=1.04100208920608; "N"=1.012756881572308 }
    21 = @{ "B"=1.02; "C"=1.018833102242169; "D"=1.028727115824554; "E"=1.020109110339748; "F"=1.036337205173399; "I"=1.026762812919656; "J"=1.025544952832684; "K"=1.032334532693093; "L"=1.023749409594542; "M"=1.039916088806156; "N"=1.012455395783456 }
    22 = @{ "B"=1.02; "C"=1.018024013616803; "D"=1.027962693101828; "E"=1.019420094225628; "F"=1.035521586759796; "I"=1.026767865498438; "J"=1.024991177845222; "K"=1.031703144176385; "L"=1.023194449864697; "M"=1.039232555522422; "N"=1.012265557879299 }
    23 = @{ "B"=1.02; "C"=1.018452901855319; "D"=1.028367887513426; "E"=1.019785297080056; "F"=1.03595397261026; "I"=1.026765424811429; "J"=1.025284767350189; "K"=1.032037873767386; "L"=1.023488648617408; "M"=1.039594977886888; "N"=1.012366221264578 }
    24 = @{ "B"=1.02; "C"=1.020141955507828; "D"=1.029964001796204; "E"=1.021224346104748; "F"=1.037655991185263; "I"=1.026750537745551; "J"=1.026440096397809; "K"=1.033355262872663; "L"=1.024646792652193; "M"=1.041020283918662; "N"=1.012761931245851 }
    25 = @{ "B"=1.02; "C"=1.022103213693989; "D"=1.031818086138188; "E"=1.022896940423851; "F"=1.039630592632542; "I"=1.026722295101752; "J"=1.027779758693215; "K"=1.034883182254064; "L"=1.025990580256047; "M"=1.042671131915277; "N"=1.013219904122195 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
